# Change the table style applied to the comparison table on slide 16
# from the deck's custom "Table_0" style to the PowerPoint built-in
# "Medium Style 2 - Accent 1" table style ({101DE442-C971-42F3-8BA1-4D48E4A8FC61}).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table
$tbl.ApplyStyle("{101DE442-C971-42F3-8BA1-4D48E4A8FC61}")
